$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D6 changes from "영역" to "사용가능 범위"
$ws.Range("D6").Value = "사용가능 범위"

# Update the saved selection to match the authored state (G17)
$ws.Range("G17").Select()
